$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new job posting row (JD_002) below the existing JD_001 row.
$ws.Range("A3").Value = "JD_002"
$ws.Range("B3").Value = "Senior Dotnet Engineer"
$ws.Range("C3").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 5

# The multi-line description would otherwise leave the row with an explicit
# custom height; auto-fitting restores the sheet's default row height
# behaviour (matching the plain data rows above it).
$ws.Rows.Item(3).AutoFit()

$wb.Save()
